# The plugin now returns an array of objects. Each object represents a row.
# Several cells that used to hold a value are now missing (sparse rows),
# and the last row's D value becomes the literal string "ende" instead of
# a number.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("First")

# --- "First" sheet: punch holes into the previously-dense grid -------------
$ws1.Range("A3").ClearContents()
$ws1.Range("D5").ClearContents()
$ws1.Range("B6").ClearContents()
$ws1.Range("D9").ClearContents()
$ws1.Range("C10").ClearContents()
$ws1.Range("B11").ClearContents()
$ws1.Range("B13").ClearContents()

# Last row's D value becomes a text marker instead of a number.
$ws1.Range("D14").Value = "ende"

# --- Selection / active sheet ----------------------------------------------
# "Second" was previously the active tab with B11 cleared last; now "First"
# is the active tab and the selection sits on the just-cleared B11 cell.
$ws1.Range("B11").Select()
